# Tutorial 6 solution update: switch date separators from "/" to "-"
# and update the Total/Real/Invalid/Absent counters for several rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextDate($cellAddress, $text) {
    # Force the cell to be treated as plain text so Excel does not
    # reinterpret the dd-mm-yyyy string as a date value, then clear the
    # temporary "@" number format so the cell keeps its original
    # (unstyled) appearance.
    $rng = $ws.Range($cellAddress)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

# --- Update the date strings in column A (rows 3-21) ---
Set-TextDate "A3"  "28-07-2022"
Set-TextDate "A4"  "01-08-2022"
Set-TextDate "A5"  "04-08-2022"
Set-TextDate "A6"  "08-08-2022"
Set-TextDate "A7"  "11-08-2022"
Set-TextDate "A8"  "15-08-2022"
Set-TextDate "A9"  "18-08-2022"
Set-TextDate "A10" "22-08-2022"
Set-TextDate "A11" "25-08-2022"
Set-TextDate "A12" "29-08-2022"
Set-TextDate "A13" "01-09-2022"
Set-TextDate "A14" "05-09-2022"
Set-TextDate "A15" "08-09-2022"
Set-TextDate "A16" "12-09-2022"
Set-TextDate "A17" "15-09-2022"
Set-TextDate "A18" "19-09-2022"
Set-TextDate "A19" "22-09-2022"
Set-TextDate "A20" "26-09-2022"
Set-TextDate "A21" "29-09-2022"

# --- Update the attendance counters (D/E/F/G/H) for the affected rows ---
# Row 3: Total Attendance -> 1, Invalid -> 1 (Absent stays 1)
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

# Row 5: Total Attendance -> 1, Real -> 1, Absent -> 0
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("H5").Value = 0

# Row 6: Total Attendance -> 1, Real -> 1, Absent -> 0
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("H6").Value = 0

# Row 11: Total Attendance -> 1, Real -> 1, Absent -> 0
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1
$ws.Range("H11").Value = 0

# Row 12: Total Attendance -> 1, Real -> 1, Absent -> 0
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1
$ws.Range("H12").Value = 0

# Row 13: Total Attendance -> 1, Real -> 1, Absent -> 0
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("H13").Value = 0

# Row 15: Total Attendance -> 1, Real -> 1, Absent -> 0
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 1
$ws.Range("H15").Value = 0

# Row 16: Total Attendance -> 1, Real -> 1, Absent -> 0
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 1
$ws.Range("H16").Value = 0

# Row 20: Total Attendance -> 1, Invalid -> 1 (Absent stays 1)
$ws.Range("D20").Value = 1
$ws.Range("G20").Value = 1
